# Update cryptocurrency price (column D) and hourly volume-change (column E)
# figures in the "cryptos" worksheet, per the GitHub Actions scheduled refresh.
#
# Column D values are stored as plain text in the workbook (e.g. "27.754.68"
# uses dots as thousands separators, which is not a valid number), so when a
# replacement value would otherwise be auto-recognised by Excel as a genuine
# number we prefix it with a leading apostrophe to force text entry and then
# immediately restore the "Normal" cell style so the quote-prefix formatting
# does not stick to the cell (matching the original unstyled text cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.751.90"
$ws.Cells.Item(2, 5).Value = "  +0.51%  "
$ws.Cells.Item(3, 4).Value = "1.850.48"
$ws.Cells.Item(3, 5).Value = "  +0.01%  "
$ws.Cells.Item(4, 5).Value = "  -1.29%  "
$ws.Cells.Item(5, 4).Value = "'318.63"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.94%  "
$ws.Cells.Item(6, 5).Value = "  -1.35%  "
$ws.Cells.Item(7, 5).Value = "  -1.57%  "
$ws.Cells.Item(8, 4).Value = "'0.3749"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.09%  "
$ws.Cells.Item(9, 4).Value = "'0.07344"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.68%  "
$ws.Cells.Item(10, 4).Value = "'0.8757"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.69%  "
$ws.Cells.Item(11, 4).Value = "'21.52"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.01%  "
$ws.Cells.Item(12, 4).Value = "1.851.24"
$ws.Cells.Item(12, 5).Value = "  -0.12%  "
$ws.Cells.Item(13, 4).Value = "'6.744"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.47%  "
$ws.Cells.Item(14, 4).Value = "'5.441"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.21%  "
$ws.Cells.Item(15, 4).Value = "'0.07140"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.03%  "
$ws.Cells.Item(16, 4).Value = "'88.88"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +4.59%  "
$ws.Cells.Item(17, 4).Value = "'1.014"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.66%  "
$ws.Cells.Item(18, 4).Value = "'0.000009001"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.90%  "
$ws.Cells.Item(19, 5).Value = "  -1.31%  "
$ws.Cells.Item(20, 4).Value = "'15.43"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.19%  "
$ws.Cells.Item(21, 4).Value = "27.765.08"
$ws.Cells.Item(21, 5).Value = "  +0.49%  "
$ws.Cells.Item(22, 4).Value = "'5.215"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.33%  "
$ws.Cells.Item(23, 5).Value = "  -1.50%  "
$ws.Cells.Item(24, 4).Value = "2.078.16"
$ws.Cells.Item(24, 5).Value = "  -0.43%  "
$ws.Cells.Item(25, 4).Value = "'1.983"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.10%  "
$ws.Cells.Item(26, 4).Value = "'155.38"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.21%  "
$ws.Cells.Item(27, 5).Value = "  -0.90%  "
$ws.Cells.Item(28, 4).Value = "'2.179"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +9.91%  "
$ws.Cells.Item(29, 5).Value = "  +0.25%  "
$ws.Cells.Item(30, 4).Value = "'118.96"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.25%  "
$ws.Cells.Item(31, 4).Value = "'0.08939"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.61%  "
$ws.Cells.Item(32, 4).Value = "'1.227"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.35%  "
$ws.Cells.Item(33, 4).Value = "'0.7776"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.34%  "
$ws.Cells.Item(34, 4).Value = "'4.540"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.27%  "
$ws.Cells.Item(35, 4).Value = "'2.925"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -2.24%  "
$ws.Cells.Item(36, 5).Value = "  -1.40%  "
$ws.Cells.Item(37, 4).Value = "'1.131"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.68%  "
$ws.Cells.Item(38, 5).Value = "  +0.63%  "
$ws.Cells.Item(39, 4).Value = "'0.05340"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.62%  "
$ws.Cells.Item(40, 4).Value = "'2.898"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.56%  "
$ws.Cells.Item(41, 4).Value = "'7.150"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +4.46%  "
$ws.Cells.Item(42, 5).Value = "  +1.13%  "
$ws.Cells.Item(43, 4).Value = "'0.5134"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.00%  "
$ws.Cells.Item(44, 4).Value = "'8.817"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.03%  "
$ws.Cells.Item(45, 4).Value = "'10.67"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.56%  "
$ws.Cells.Item(46, 4).Value = "'107.32"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.39%  "
$ws.Cells.Item(47, 5).Value = "  +1.38%  "
$ws.Cells.Item(48, 4).Value = "'0.06468"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.81%  "
$ws.Cells.Item(49, 4).Value = "'1.012"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.55%  "
$ws.Cells.Item(50, 4).Value = "'1.690"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.73%  "
$ws.Cells.Item(51, 5).Value = "  -2.66%  "
